# Fruta / hortaliza, semanal
# Inserts two new weekly price rows (Primera / Segunda) right before the
# existing row 50 (pushing the rest of the block down by two rows), and
# appends one more new weekly price row (Primera) after the last existing
# row (which ends up at row 70 once shifted).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two blank rows above row 50 (existing data shifts down) -------
$ws.Rows.Item(50).Resize(2).EntireRow.Insert()

function Set-PapayaRow($row, $fecha, $calidad, $volumen, $precio, $precioKg, $kgUnidad) {
    $ws.Cells.Item($row, 1).Value = 3
    $ws.Cells.Item($row, 2).Value = "Femacal de La Calera"
    $ws.Cells.Item($row, 3).Value = "Coquimbo"
    $ws.Cells.Item($row, 4).Value = $fecha
    $ws.Cells.Item($row, 5).Value = 5
    $ws.Cells.Item($row, 6).Value = "Fruta"
    $ws.Cells.Item($row, 7).Value = 100108
    $ws.Cells.Item($row, 8).Value = "Tropicales y subtropicales"
    $ws.Cells.Item($row, 9).Value = 100108004
    $ws.Cells.Item($row, 10).Value = "Papaya"
    $ws.Cells.Item($row, 11).Value = "Cultivar IV Región"
    $ws.Cells.Item($row, 12).Value = $calidad
    $ws.Cells.Item($row, 13).Value = $volumen
    $ws.Cells.Item($row, 14).Value = $precio
    $ws.Cells.Item($row, 15).Value = $precio
    $ws.Cells.Item($row, 16).Value = $precio
    $ws.Cells.Item($row, 17).Value = "`$/bandeja 10 kilos"
    $ws.Cells.Item($row, 18).Value = "Provincia del Elquí"
    $ws.Cells.Item($row, 19).Value = $precioKg
    $ws.Cells.Item($row, 20).Value = $kgUnidad
}

# New row 50: 2023-07-13, Primera
Set-PapayaRow 50 45120 "Primera" 56 18000 1800 10

# New row 51: 2023-07-13, Segunda
Set-PapayaRow 51 45120 "Segunda" 50 15000 1500 10

# --- Append one more new row after the (now shifted) last existing row ----
# The previously-last data row (68) is now row 70, so the new row goes at 71.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
$newRow = $lastRow + 1

Set-PapayaRow $newRow 45121 "Primera" 56 20000 2000 10

# The appended row is not produced via a row-insert, so it does not inherit
# the date column's number format automatically - copy it explicitly so the
# new date cell renders/serialises the same way as the rest of column D.
$ws.Cells.Item($newRow, 4).NumberFormat = $ws.Cells.Item($newRow - 1, 4).NumberFormat
